$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.12926056049198564
$ws.Range("B1").Value = 0.12910707082606621
$ws.Range("A2").Value = -0.082999324368310923
$ws.Range("B2").Value = 0.082533622791968853
$ws.Range("A3").Value = 0.079113465637554725
$ws.Range("B3").Value = -0.079378614891993493
$ws.Range("A4").Value = -0.20461053671048646
$ws.Range("B4").Value = 0.203508194453633
$ws.Range("A5").Value = -0.19750819526314345
$ws.Range("B5").Value = 0.1952776495935451
$ws.Range("A6").Value = -0.10938719146744047
$ws.Range("B6").Value = 0.1092071166594577
$ws.Range("A7").Value = -0.089207117643365308
$ws.Range("B7").Value = 0.088753656558134608
$ws.Range("A8").Value = -0.068753657553519254
$ws.Range("B8").Value = 0.068361062918652848
$ws.Range("A9").Value = -0.062361063790210558
$ws.Range("B9").Value = 0.0620268134949713
$ws.Range("A10").Value = -0.056026814380217616
$ws.Range("B10").Value = 0.055983089476733028
$ws.Range("A11").Value = -0.051483090349066174
$ws.Range("B11").Value = 0.051404215030348865
$ws.Range("A12").Value = -0.045404215920569424
$ws.Range("B12").Value = 0.045153683172783232
$ws.Range("A13").Value = -0.039153684076731032
$ws.Range("B13").Value = 0.039085762061716878
$ws.Range("A14").Value = -0.027085763028678933
$ws.Range("B14").Value = 0.02705335577548329
$ws.Range("A15").Value = -0.021053356686577374
$ws.Range("B15").Value = 0.021027906923847262
$ws.Range("A16").Value = -0.015027907837892762
$ws.Range("B16").Value = 0.015004110355869349
$ws.Range("A17").Value = -0.009004111273759996
$ws.Range("B17").Value = 0.0089999990517437567
$ws.Range("A18").Value = -0.036107747289221948
$ws.Range("B18").Value = 0.036095983882642457
$ws.Range("A19").Value = -0.027095984704586407
$ws.Range("B19").Value = 0.027013029824263679
$ws.Range("A20").Value = -0.018013030653115436
$ws.Range("B20").Value = 0.018004186997368166
$ws.Range("A21").Value = -0.0090041878272053566
$ws.Range("B21").Value = 0.0089999991694398318
$ws.Range("A22").Value = -0.093934217152774835
$ws.Range("B22").Value = 0.093626084992985525
$ws.Range("A23").Value = -0.084626085829043518
$ws.Range("B23").Value = 0.084124969565106511
$ws.Range("A24").Value = -0.042124970738264267
$ws.Range("B24").Value = 0.041999998820732465
$ws.Range("A25").Value = -0.0798904589683751
$ws.Range("B25").Value = 0.079812624263482235
$ws.Range("A26").Value = -0.073812625108391927
$ws.Range("B26").Value = 0.073715863077918442
$ws.Range("A27").Value = -0.072805871519498133
$ws.Range("B27").Value = 0.072405883253004344
$ws.Range("A28").Value = -0.06640588411642856
$ws.Range("B28").Value = 0.066151414040763967
$ws.Range("A29").Value = -0.054151414973283352
$ws.Range("B29").Value = 0.054049196513540565
$ws.Range("A30").Value = -0.034049197529342212
$ws.Range("B30").Value = 0.033820579220295244
$ws.Range("A31").Value = -0.027017915712535867
$ws.Range("B31").Value = 0.027000556983464818
$ws.Range("A32").Value = -0.0060005580171358375
$ws.Range("B32").Value = 0.0059999991135226693

$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334
